$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. O1468 changes from 0 to 1 ---
$ws.Cells.Item(1468, 15).Value = 1

# --- 2. R1470 and R1471 change from empty/inlineStr to numeric 0 ---
$ws.Cells.Item(1470, 18).Value = 0
$ws.Cells.Item(1471, 18).Value = 0

$newRows = @(
    ,@(1472, 45474, 614.2000122070312, 622.9000244140625, 601, 615.3499755859375, 615.3499755859375, 14670413, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0)
    ,@(1473, 45481, 632, 655.7999877929688, 629, 650.0999755859375, 650.0999755859375, 21739511, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0)
    ,@(1474, 45488, 658.9500122070312, 686.25, 646.7000122070312, 668.6500244140625, 668.6500244140625, 12824804, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 2)
    ,@(1475, 45495, 669.9500122070312, 682, 649.0999755859375, 679.75, 679.75, 10953948, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0)
    ,@(1476, 45502, 679.5, 691, 659.5999755859375, 662.4000244140625, 662.4000244140625, 8734092, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0)
    ,@(1477, 45509, 656.1500244140625, 682, 626.25, 653, 653, 23600936, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0)
    ,@(1478, 45516, 648, 665.75, 637.7000122070312, 661.0499877929688, 661.0499877929688, 8854407, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0)
    ,@(1479, 45523, 666.2000122070312, 688.6500244140625, 660.7999877929688, 678.2000122070312, 678.2000122070312, 7380555, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0)
    ,@(1480, 45530, 681, 692.7999877929688, 644.75, 647.1500244140625, 647.1500244140625, 15620551, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0)
    ,@(1481, 45537, 648, 673.9500122070312, 633, 665.25, 665.25, 16463119, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 0)
    ,@(1482, 45544, 665, 690.2000122070312, 660.5, 681.9500122070312, 681.9500122070312, 11938554, 2024, 9, 9, 0, 0, 0, 37, 0, 0, 0)
    ,@(1483, 45551, 683.0499877929688, 710.5, 675.5499877929688, 709, 709, 9077184, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0)
    ,@(1484, 45558, 709.0499877929688, 713.5, 681.5499877929688, 692.4500122070312, 692.4500122070312, 21963579, 2024, 9, 23, 0, 0, 0, 39, 0, 0, 0)
    ,@(1485, 45565, 680.5, 719.8499755859375, 680.5, 690.2000122070312, 690.2000122070312, 11018938, 2024, 9, 30, 0, 0, 0, 40, 1, 0, 0)
    ,@(1486, 45572, 686.25, 704.6500244140625, 673.7999877929688, 685.5, 685.5, 7005441, 2024, 10, 7, 0, 0, 0, 41, 0, 0, 0)
    ,@(1487, 45579, 685.5, 692.5, 654.5999755859375, 669.2999877929688, 669.2999877929688, 9593468, 2024, 10, 14, 0, 0, 0, 42, 0, 0, 0)
    ,@(1488, 45586, 665.0999755859375, 670.0999755859375, 631.0499877929688, 640.0999755859375, 640.0999755859375, 7542971, 2024, 10, 21, 0, 0, 0, 43, 0, 0, 0)
    ,@(1489, 45593, 640.0999755859375, 687, 616.8499755859375, 645.9500122070312, 645.9500122070312, 21370838, 2024, 10, 28, 0, 0, 0, 44, 0, 0, 0)
    ,@(1490, 45600, 642.5999755859375, 650, 624.7999877929688, 629.8499755859375, 629.8499755859375, 6352800, 2024, 11, 4, 0, 0, 0, 45, 0, 0, 0)
)
# --- 3. Append 19 new weekly rows (1472-1490) ---
foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
    $ws.Cells.Item($r, 12).Value = $row[12]
    $ws.Cells.Item($r, 13).Value = $row[13]
    $ws.Cells.Item($r, 14).Value = $row[14]
    $ws.Cells.Item($r, 15).Value = $row[15]
    $ws.Cells.Item($r, 16).Value = $row[16]
    $ws.Cells.Item($r, 17).Value = $row[17]
}
